# Generate Report for Handoff
# Updates the localization status report: the handoff status moves from
# "Ready for handoff" to "In Translation", and the associated timestamps
# advance to the moment the new handoff report was generated.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: per-language status + latest handoff-xliff generation date
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("G2").Value = "2017-03-02 06:39:45"

# zh-cn detail sheet: status + latest handoff datetime
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("H2").Value = "2017-03-02 06:39:29"

# de-de detail sheet: status + latest handoff datetime
$dede.Range("C2").Value = "In Translation"
$dede.Range("H2").Value = "2017-03-02 06:39:45"

# The "Status" column(s) got a little narrower now that the longest value is
# "In Translation" rather than "Ready for handoff" - shrink them to match
# (Excel snaps ColumnWidth to its internal pixel grid, so this lands as close
# to the regenerated report's width as that grid allows).
$overview.Range("E1").EntireColumn.ColumnWidth = 12.52
$overview.Range("F1").EntireColumn.ColumnWidth = 12.52
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.52
$dede.Range("C1").EntireColumn.ColumnWidth = 12.52
